# multi browser implementation for chrome and firefox and accounting cash
# -------------------------------------------------------------------------
# Repayment Schedule recompute: the "Outstanding" (O) column is cleared out
# (it was redundant with the final "Outstanding" figure already carried in
# column P), and the interest/principal/balance figures for the later
# installments are refreshed to the newly recalculated amounts. The active
# sheet/selection also moves from "Repayment Schedule" back to
# "NewLoanInput".

$wb = $excel.ActiveWorkbook

$wsInput    = $wb.Worksheets.Item("NewLoanInput")
$wsSchedule = $wb.Worksheets.Item("Repayment Schedule")

# --- Repayment Schedule: drop the stray "Outstanding" (O) column values ---
# Row 2's empty placeholder cell shifts from P2 over to O2 (same blank style)
$wsSchedule.Range("N2").Copy()
$wsSchedule.Range("O2").PasteSpecial(-4122)
$wsSchedule.Range("P2").Clear()

# Rows 3-14 simply drop the redundant zero in the O column entirely
$wsSchedule.Range("O3:O14").Clear()

# --- Repayment Schedule: refresh recalculated figures for installments 4-12 ---
$wsSchedule.Range("F6").Value  = 812.14
$wsSchedule.Range("G6").Value  = 6822.39
$wsSchedule.Range("H6").Value  = 76.349999999999994

$wsSchedule.Range("F7").Value  = 820.27
$wsSchedule.Range("G7").Value  = 6002.12
$wsSchedule.Range("H7").Value  = 68.22

$wsSchedule.Range("F8").Value  = 828.47
$wsSchedule.Range("G8").Value  = 5173.6499999999996
$wsSchedule.Range("H8").Value  = 60.02

$wsSchedule.Range("F9").Value  = 836.75
$wsSchedule.Range("G9").Value  = 4336.8999999999996
$wsSchedule.Range("H9").Value  = 51.74

$wsSchedule.Range("F10").Value = 845.12
$wsSchedule.Range("G10").Value = 3491.78
$wsSchedule.Range("H10").Value = 43.37

$wsSchedule.Range("F11").Value = 853.57
$wsSchedule.Range("G11").Value = 2638.21
$wsSchedule.Range("H11").Value = 34.92

$wsSchedule.Range("F12").Value = 862.11
$wsSchedule.Range("G12").Value = 1776.1
$wsSchedule.Range("H12").Value = 26.38

$wsSchedule.Range("F13").Value = 870.73
$wsSchedule.Range("G13").Value = 905.37
$wsSchedule.Range("H13").Value = 17.760000000000002

$wsSchedule.Range("F14").Value = 905.37
$wsSchedule.Range("H14").Value = 9.0500000000000007
$wsSchedule.Range("K14").Value = 914.42
$wsSchedule.Range("P14").Value = 914.42

# --- View state: active sheet moves from "Repayment Schedule" back to ---
# --- "NewLoanInput", with a fresh selection on each sheet.            ---
[void]$wsSchedule.Range("A2:XFD19").Select()

$wsInput.Activate()
[void]$wsInput.Range("B24").Select()
